$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '41.598.30'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.39%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.164.86'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.28%  '

$ws.Range("E4").Value = '  -0.04%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '237.76'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.39%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.606'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -3.02%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '71.04'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.63%  '

$ws.Range("E8").Value = '  +0.01%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.573'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.37%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '39.67'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.45%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0900'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -4.80%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.04'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.47%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0998'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.67%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.67'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.60%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.489.32'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.41%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '14.23'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.165.80'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.85%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.783'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -6.03%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '41.464.73'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.47%  '

$ws.Range("E20").Value = '  -3.71%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '69.58'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.75%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.74'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -6.69%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.02'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -9.28%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '226.65'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.03%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.25%  '

$ws.Range("E26").Value = '  -0.11%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.67'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -6.13%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.29'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -9.27%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -3.67%  '

$ws.Range("E30").Value = '  -1.08%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '171.27'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.36%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '19.73'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.17%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '33.42'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +11.71%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0773'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.81%  '

$ws.Range("E35").Value = '  -6.79%  '

$ws.Range("E36").Value = '  -3.58%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.105'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.91%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.26'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.07%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0301'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.23%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '12.14'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -8.86%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.06'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.41%  '

$ws.Range("E42").Value = '  -4.81%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '58.62'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -8.04%  '

$ws.Range("E44").Value = '  -3.84%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.36'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.85%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0960'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.86%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '95.60'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -6.86%  '

$ws.Range("E48").Value = '  -2.53%  '

$ws.Range("E49").Value = '  -4.70%  '

$ws.Range("E50").Value = '  -7.30%  '

$ws.Range("E51").Value = '  -2.40%  '

